$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of B9:D36 (values + types), but keep cell styles/formatting.
$ws.Range("B9:D36").ClearContents()

# Update the active selection to D12 as per the diff.
$ws.Range("D12").Select()
